$wb = $excel.ActiveWorkbook

# --- sheet1 ("user-data"): add row 4 (ReaderID 3 / bibi / 123 / bear) ---
$ws1 = $wb.Worksheets.Item("user-data")
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "bibi"
# "123" must stay text (matches existing Password column, which is
# numberStoredAsText) - a leading apostrophe forces text storage instead
# of Excel re-interpreting the digits as a number.
$ws1.Range("C4").Value = "'123"
$ws1.Range("D4").Value = "bear"

# --- sheet2 ("book-data"): add row 7 (ReaderID 3 / In Cold Blood / Truman Capote ) ---
$ws2 = $wb.Worksheets.Item("book-data")
# column A on this sheet is entirely text (numberStoredAsText), so force "3"
# to stay text rather than become the number 3.
$ws2.Range("A7").Value = "'3"
$ws2.Range("B7").Value = "In Cold Blood"
$ws2.Range("C7").Value = "Truman Capote "
